$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.344.96"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "'2.273.59"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  -0.90%  "
$ws.Range("D5").Value = "'300.65"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "'95.24"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").Value = "'0.567"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").Value = "'34.22"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("D11").Value = "'0.0795"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "'7.23"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "'2.616.89"
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "'2.270.53"
$ws.Range("E15").Value = "  -3.53%  "
$ws.Range("D16").Value = "'0.813"
$ws.Range("E16").Value = "  -1.88%  "
$ws.Range("D17").Value = "'13.66"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("D18").Value = "'45.149.67"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("D19").Value = "'13.22"
$ws.Range("E19").Value = "  +13.12%  "
$ws.Range("D20").Value = "'0.0₃0918"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "'6.04"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").Value = "'65.50"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "'240.59"
$ws.Range("E23").Value = "  +1.54%  "
$ws.Range("D24").Value = "'2.88"
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("D26").Value = "'1.92"
$ws.Range("E26").Value = "  -3.23%  "
$ws.Range("D27").Value = "'41.21"
$ws.Range("E27").Value = "  +10.91%  "
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("D29").Value = "'9.59"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").Value = "'19.68"
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").Value = "'151.71"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "'5.54"
$ws.Range("E32").Value = "  -4.73%  "
$ws.Range("D33").Value = "'0.0792"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  -1.87%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "'2.92"
$ws.Range("E36").Value = "  -7.06%  "
$ws.Range("D37").Value = "'0.104"
$ws.Range("E37").Value = "  -3.18%  "
$ws.Range("E38").Value = "  -5.67%  "
$ws.Range("D39").Value = "'3.92"
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("E40").Value = "  +5.55%  "
$ws.Range("D41").Value = "'3.23"
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").Value = "'13.52"
$ws.Range("E42").Value = "  -8.94%  "
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("E44").Value = "  +11.22%  "
$ws.Range("D45").Value = "'1.767.02"
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("D47").Value = "'70.05"
$ws.Range("E47").Value = "  +2.54%  "
$ws.Range("D48").Value = "'76.46"
$ws.Range("E48").Value = "  -4.59%  "
$ws.Range("D49").Value = "'95.59"
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("D50").Value = "'4.74"
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").Value = "'53.65"
$ws.Range("E51").Value = "  +0.17%  "
